$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. The data table is being extended one more year (2023) into a new
#    column K, mirroring the existing year columns (B..J = 2014..2022).
# ------------------------------------------------------------------

# Copy the formatting (font, borders, number format, alignment) of the
# last existing year column (J) into the new column K for every row of
# the little table, then write the 2023 figures into it.
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3:K6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 10761
$ws.Range("K5").Value = 3005
$ws.Range("K6").Value = 7756

# ------------------------------------------------------------------
# 2. Border touch-up: the header row (3, the year row) used to be
#    boxed on both the top and the bottom; now that the table keeps
#    growing to the right it only keeps the top rule, and the new
#    rightmost column (K) gets a closing right-hand rule on every row
#    of the block so the table looks closed again.
# ------------------------------------------------------------------
$ws.Range("B3:K3").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none
$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1        # xlEdgeRight -> thin
$ws.Range("K3:K6").Borders.Item(10).Weight = 2            # xlThin

# ------------------------------------------------------------------
# 3. Pre-widen the columns to the right of the table (K..S) the same
#    way B..J already are, so future year columns can be dropped in
#    without reformatting again.
# ------------------------------------------------------------------
$ws.Columns("K:S").ColumnWidth = $ws.Columns("J").ColumnWidth
